$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.486.99"
$ws.Range("E2").Value = "  -1.55%  "
$ws.Range("D3").Value = "2.901.06"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.42%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.549"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.14%  "
$ws.Range("D9").Value = "2.905.70"
$ws.Range("E9").Value = "  -2.50%  "
$ws.Range("E10").Value = "  -5.19%  "
$ws.Range("E11").Value = "  -2.74%  "
$ws.Range("E12").Value = "  -2.42%  "
$ws.Range("D13").Value = "3.409.99"
$ws.Range("E13").Value = "  -2.33%  "
$ws.Range("E14").Value = "  +2.39%  "
$ws.Range("D15").Value = "60.507.39"
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.97%  "
$ws.Range("D17").Value = "2.901.28"
$ws.Range("E17").Value = "  -2.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000140"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.50%  "
$ws.Range("E19").Value = "  -3.60%  "
$ws.Range("E20").Value = "  -3.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.452"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.64%  "
$ws.Range("E27").Value = "  -5.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.80%  "
$ws.Range("D30").Value = "0.0₃0846"
$ws.Range("E30").Value = "  -9.55%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("E33").Value = "  -4.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.60%  "
$ws.Range("E35").Value = "  -6.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.83%  "
$ws.Range("E38").Value = "  -5.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.69"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.46"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.36%  "
$ws.Range("E41").Value = "  -5.09%  "
$ws.Range("D42").Value = "2.286.45"
$ws.Range("E42").Value = "  -4.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.647"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0580"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.66%  "
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.49%  "
$ws.Range("E48").Value = "  -3.40%  "
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("E50").Value = "  -2.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "248.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.68%  "
